# Remove the "total=N" tag values that were appended to the *_tags columns
# of the measure-score example sheets (IAR-DST, K10+, K5, SDQ). These were
# redundant "total=" annotations baked into the tags column of each scored
# row; the commit removes them (leaving the other tag values, if any,
# untouched) and lets the shared-string table shrink accordingly.

$wb = $excel.ActiveWorkbook

$sheetNames = @("IAR-DST", "K10+", "K5", "SDQ")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string] -and $val.StartsWith("total=")) {
                $cell.ClearContents()
            }
        }
    }
}
